$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: a neighboring cell with the plain default style, used as a style donor
# so that writing text-looking numeric strings does not alter the cell style index.
$donorStyle = $ws.Range("B2").Style

function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $donorStyle
}

# Row 15: swap coin identity
Set-TextValue $ws.Range("B15") "WrappedEther"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D15") "1.777.53"
Set-TextValue $ws.Range("E15") "  -1.03%  "

# Row 16: swap coin identity
Set-TextValue $ws.Range("B16") "Chainlink"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D16") "7.197"
Set-TextValue $ws.Range("E16") "  -2.04%  "

# Row 42: swap coin identity
Set-TextValue $ws.Range("B42") "Aptos"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D42") "10.98"
Set-TextValue $ws.Range("E42") "  -5.44%  "

# Row 43: swap coin identity
Set-TextValue $ws.Range("B43") "TrustWalletToken"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D43") "1.171"
Set-TextValue $ws.Range("E43") "  +1.14%  "

# Row 2
Set-TextValue $ws.Range("D2") "27.955.63"
Set-TextValue $ws.Range("E2") "  -1.05%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.776.40"
Set-TextValue $ws.Range("E3") "  -1.30%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.27%  "

# Row 5
Set-TextValue $ws.Range("D5") "315.56"
Set-TextValue $ws.Range("E5") "  -0.32%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.004"
Set-TextValue $ws.Range("E6") "  +0.27%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.5354"
Set-TextValue $ws.Range("E7") "  -2.56%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3724"
Set-TextValue $ws.Range("E8") "  -3.88%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.07398"
Set-TextValue $ws.Range("E9") "  -2.46%  "

# Row 10
Set-TextValue $ws.Range("D10") "41.47"
Set-TextValue $ws.Range("E10") "  -3.09%  "

# Row 11
Set-TextValue $ws.Range("D11") "1.090"
Set-TextValue $ws.Range("E11") "  -2.57%  "

# Row 12
Set-TextValue $ws.Range("D12") "1.004"
Set-TextValue $ws.Range("E12") "  +0.17%  "

# Row 13
Set-TextValue $ws.Range("D13") "20.42"
Set-TextValue $ws.Range("E13") "  -3.35%  "

# Row 14
Set-TextValue $ws.Range("D14") "6.064"
Set-TextValue $ws.Range("E14") "  -2.29%  "

# Row 17
Set-TextValue $ws.Range("D17") "88.10"
Set-TextValue $ws.Range("E17") "  -4.28%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.00001050"
Set-TextValue $ws.Range("E18") "  -1.99%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.06466"
Set-TextValue $ws.Range("E19") "  +0.12%  "

# Row 20
Set-TextValue $ws.Range("E20") "  +0.17%  "

# Row 21
Set-TextValue $ws.Range("D21") "17.35"
Set-TextValue $ws.Range("E21") "  +0.44%  "

# Row 22
Set-TextValue $ws.Range("D22") "5.886"
Set-TextValue $ws.Range("E22") "  -1.46%  "

# Row 23
Set-TextValue $ws.Range("D23") "28.001.51"
Set-TextValue $ws.Range("E23") "  -0.94%  "

# Row 24
Set-TextValue $ws.Range("D24") "11.07"
Set-TextValue $ws.Range("E24") "  -3.82%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.094"
Set-TextValue $ws.Range("E25") "  -3.25%  "

# Row 26
Set-TextValue $ws.Range("D26") "157.71"
Set-TextValue $ws.Range("E26") "  -0.34%  "

# Row 27
Set-TextValue $ws.Range("D27") "20.17"
Set-TextValue $ws.Range("E27") "  -2.48%  "

# Row 28
Set-TextValue $ws.Range("D28") "1.979.21"
Set-TextValue $ws.Range("E28") "  -1.36%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.276"
Set-TextValue $ws.Range("E29") "  -5.63%  "

# Row 30
Set-TextValue $ws.Range("D30") "119.81"
Set-TextValue $ws.Range("E30") "  -3.06%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.097"
Set-TextValue $ws.Range("E31") "  -3.15%  "

# Row 32
Set-TextValue $ws.Range("E32") "  +2.17%  "

# Row 33
Set-TextValue $ws.Range("D33") "3.652"
Set-TextValue $ws.Range("E33") "  -0.56%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.477"
Set-TextValue $ws.Range("E34") "  -4.44%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.2231"
Set-TextValue $ws.Range("E35") "  -4.41%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.06356"
Set-TextValue $ws.Range("E36") "  -0.39%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.02260"
Set-TextValue $ws.Range("E37") "  -2.56%  "

# Row 38
Set-TextValue $ws.Range("D38") "4.959"
Set-TextValue $ws.Range("E38") "  -1.29%  "

# Row 39
Set-TextValue $ws.Range("D39") "8.424"
Set-TextValue $ws.Range("E39") "  -5.28%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.6145"
Set-TextValue $ws.Range("E40") "  -3.98%  "

# Row 41
Set-TextValue $ws.Range("E41") "  +3.66%  "

# Row 44
Set-TextValue $ws.Range("D44") "1.003"
Set-TextValue $ws.Range("E44") "  +0.29%  "

# Row 45
Set-TextValue $ws.Range("D45") "13.30"
Set-TextValue $ws.Range("E45") "  -1.04%  "

# Row 46
Set-TextValue $ws.Range("D46") "3.663"
Set-TextValue $ws.Range("E46") "  -0.52%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.5740"
Set-TextValue $ws.Range("E47") "  -4.18%  "

# Row 48
Set-TextValue $ws.Range("D48") "125.62"
Set-TextValue $ws.Range("E48") "  +1.34%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.194"

# Row 50
Set-TextValue $ws.Range("D50") "1.921"
Set-TextValue $ws.Range("E50") "  -2.73%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.06817"
Set-TextValue $ws.Range("E51") "  -1.07%  "
